$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B1 holds the shared string "value" -> rename to "first_release_value".
# (A1 keeps referencing "date"; only the B-column label text changes.)
$ws.Range("B1").Value = "first_release_value"

# Extend the date-column number format/style down into the newly-added rows
# (54-84) by copying the format from the last previously-existing date cell.
$ws.Range("A53").Copy() | Out-Null
$ws.Range("A54:A84").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rewrite the full data block (dates in A, qoq values in B) for rows 2-84.
$ws.Range("A2").Value = 38398
$ws.Range("B2").Value = -0.3055163919209463
$ws.Range("A3").Value = 38487
$ws.Range("B3").Value = -0.05930270183239372
$ws.Range("A4").Value = 38579
$ws.Range("B4").Value = 0.6231390923824733
$ws.Range("A5").Value = 38671
$ws.Range("B5").Value = -0.7175835844218028
$ws.Range("A6").Value = 38763
$ws.Range("B6").Value = 0.8366283085903774
$ws.Range("A7").Value = 38852
$ws.Range("B7").Value = 0.1767276857804774
$ws.Range("A8").Value = 38944
$ws.Range("B8").Value = 0.3528382780321806
$ws.Range("A9").Value = 39036
$ws.Range("B9").Value = 0.8301537441992792
$ws.Range("A10").Value = 39128
$ws.Range("B10").Value = -1.966285411910945
$ws.Range("A11").Value = 39217
$ws.Range("B11").Value = 0.7805496083026924
$ws.Range("A12").Value = 39309
$ws.Range("B12").Value = 0.5392191980432131
$ws.Range("A13").Value = 39401
$ws.Range("B13").Value = -0.7807892990731773
$ws.Range("A14").Value = 39493
$ws.Range("B14").Value = 0.3147656593484953
$ws.Range("A15").Value = 39583
$ws.Range("B15").Value = -0.6071410908285912
$ws.Range("A16").Value = 39675
$ws.Range("B16").Value = 0.2561593687207875
$ws.Range("A17").Value = 39767
$ws.Range("B17").Value = -0.2643890380460761
$ws.Range("A18").Value = 39859
$ws.Range("B18").Value = 0.4614603479951001
$ws.Range("A19").Value = 39948
$ws.Range("B19").Value = 0.5995551687457663
$ws.Range("A20").Value = 40040
$ws.Range("B20").Value = -0.88436028068827
$ws.Range("A21").Value = 40132
$ws.Range("B21").Value = -0.1757398580474785
$ws.Range("A22").Value = 40224
$ws.Range("B22").Value = -0.8215943001740271
$ws.Range("A23").Value = 40313
$ws.Range("B23").Value = 0.6
$ws.Range("A24").Value = 40405
$ws.Range("B24").Value = 0.4
$ws.Range("A25").Value = 40497
$ws.Range("B25").Value = 0.6349193672116513
$ws.Range("A26").Value = 40589
$ws.Range("B26").Value = 0.4
$ws.Range("A27").Value = 40678
$ws.Range("B27").Value = -0.5972483967898228
$ws.Range("A28").Value = 40770
$ws.Range("B28").Value = 0.8140368792747665
$ws.Range("A29").Value = 40862
$ws.Range("B29").Value = -0.2012463990342326
$ws.Range("A30").Value = 40954
$ws.Range("B30").Value = 0.4
$ws.Range("A31").Value = 41044
$ws.Range("B31").Value = 0.1
$ws.Range("A32").Value = 41136
$ws.Range("B32").Value = 0.3
$ws.Range("A33").Value = 41228
$ws.Range("B33").Value = -0.340329042122363
$ws.Range("A34").Value = 41320
$ws.Range("B34").Value = 0.8
$ws.Range("A35").Value = 41409
$ws.Range("B35").Value = 0.6
$ws.Range("A36").Value = 41501
$ws.Range("B36").Value = 0.1
$ws.Range("A37").Value = 41593
$ws.Range("B37").Value = -0.2963573102740611
$ws.Range("A38").Value = 41685
$ws.Range("B38").Value = 0.6873521384730878
$ws.Range("A39").Value = 41774
$ws.Range("B39").Value = 0.1052962261794335
$ws.Range("A40").Value = 41866
$ws.Range("B40").Value = 0.7267133658511682
$ws.Range("A41").Value = 41958
$ws.Range("B41").Value = 0.7028875639548886
$ws.Range("A42").Value = 42050
$ws.Range("B42").Value = 0.6036537137213145
$ws.Range("A43").Value = 42139
$ws.Range("B43").Value = 0.1326972526782129
$ws.Range("A44").Value = 42231
$ws.Range("B44").Value = 0.5774359918206358
$ws.Range("A45").Value = 42323
$ws.Range("B45").Value = 0.4228511177569345
$ws.Range("A46").Value = 42415
$ws.Range("B46").Value = 0.421075791701611
$ws.Range("A47").Value = 42505
$ws.Range("B47").Value = 0.1946648784293643
$ws.Range("A48").Value = 42597
$ws.Range("B48").Value = 0.3978154615661396
$ws.Range("A49").Value = 42689
$ws.Range("B49").Value = 0.1561519231779869
$ws.Range("A50").Value = 42781
$ws.Range("B50").Value = 0.3026479405721147
$ws.Range("A51").Value = 42870
$ws.Range("B51").Value = 0.9379874529000176
$ws.Range("A52").Value = 42962
$ws.Range("B52").Value = -0.1353320131999567
$ws.Range("A53").Value = 43054
$ws.Range("B53").Value = 0.009046673504869318
$ws.Range("A54").Value = 43146
$ws.Range("B54").Value = 0.4431994051349051
$ws.Range("A55").Value = 43235
$ws.Range("B55").Value = 0.305320875257479
$ws.Range("A56").Value = 43327
$ws.Range("B56").Value = -0.3
$ws.Range("A57").Value = 43419
$ws.Range("B57").Value = 0.2158019844203096
$ws.Range("A58").Value = 43511
$ws.Range("B58").Value = 1.2
$ws.Range("A59").Value = 43600
$ws.Range("B59").Value = 0.0941018033610419
$ws.Range("A60").Value = 43692
$ws.Range("B60").Value = 0.4323758554038761
$ws.Range("A61").Value = 43784
$ws.Range("B61").Value = 0.01870135355044056
$ws.Range("A62").Value = 43876
$ws.Range("B62").Value = -3.2
$ws.Range("A63").Value = 43966
$ws.Range("B63").Value = -10.91706676560194
$ws.Range("A64").Value = 44058
$ws.Range("B64").Value = 7.25
$ws.Range("A65").Value = 44150
$ws.Range("B65").Value = -0.7999999999999972
$ws.Range("A66").Value = 44242
$ws.Range("B66").Value = -2.400000000000006
$ws.Range("A67").Value = 44331
$ws.Range("B67").Value = 3.480874220397794
$ws.Range("A68").Value = 44423
$ws.Range("B68").Value = 6.400000000000006
$ws.Range("A69").Value = 44515
$ws.Range("B69").Value = -1.099999999999994
$ws.Range("A70").Value = 44607
$ws.Range("B70").Value = 0.4209467346675666
$ws.Range("A71").Value = 44696
$ws.Range("B71").Value = 1.540000000000006
$ws.Range("A72").Value = 44788
$ws.Range("B72").Value = -0.1490000000000009
$ws.Range("A73").Value = 44880
$ws.Range("B73").Value = -0.4999966213670604
$ws.Range("A74").Value = 44972
$ws.Range("B74").Value = -0.7399999999999949
$ws.Range("A75").Value = 45061
$ws.Range("B75").Value = -0.2510000000000048
$ws.Range("A76").Value = 45153
$ws.Range("B76").Value = 0.2000000000000028
$ws.Range("A77").Value = 45245
$ws.Range("B77").Value = 0.1200000000000045
$ws.Range("A78").Value = 45337
$ws.Range("B78").Value = 0.09999999999999432
$ws.Range("A79").Value = 45427
$ws.Range("B79").Value = 0
$ws.Range("A80").Value = 45519
$ws.Range("B80").Value = 0.09999999999999432
$ws.Range("A81").Value = 45611
$ws.Range("B81").Value = 0.3200015876295765
$ws.Range("A82").Value = 45703
$ws.Range("B82").Value = 0.09999771493470178
$ws.Range("A83").Value = 45792
$ws.Range("B83").Value = 0.2000034419242951
$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = 0.03999999999999204
